# Update sample data for document extraction rules to include additional
# fields (context_before / context_after) and extraction methods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -------------------------------------------------------
# Insert two new header columns (D: context_before, E: context_after)
# and push the old "instructions" header out to column F.
$ws.Range("F1").Value = "instructions"
$ws.Range("D1").Value = "context_before"
$ws.Range("E1").Value = "context_after"
# Make sure the new header cells pick up the same bold/border/centered
# style that is already applied to the other header cells.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# --- Row 2: Invoice Number (unchanged key fields) ----------------------
$ws.Range("A2").Value = "Invoice Number"
$ws.Range("B2").Value = "Invoice #"
$ws.Range("C2").Value = "after_pattern"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "`n"
$ws.Range("F2").Value = ""

# --- Row 3: Invoice Date (was Total Amount) ----------------------------
$ws.Range("A3").Value = "Invoice Date"
$ws.Range("B3").Value = "Invoice Date:"
$ws.Range("C3").Value = "after_pattern"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "`n"
$ws.Range("F3").Value = ""

# --- Row 4: Total Amount (was Customer Name) ---------------------------
$ws.Range("A4").Value = "Total Amount"
$ws.Range("B4").Value = "Total:"
$ws.Range("C4").Value = "after_pattern"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = "`n"
$ws.Range("F4").Value = ""

# --- Row 5: Company Name (was Shipping Address) -------------------------
$ws.Range("A5").Value = "Company Name"
$ws.Range("B5").Value = "Company:"
$ws.Range("C5").Value = "after_pattern"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "`n"
$ws.Range("F5").Value = ""

# --- Row 6: Customer Address (was Order Date) ---------------------------
$ws.Range("A6").Value = "Customer Address"
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "nlp"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "Find the complete customer address in the document"

# --- Row 7: Customer Email (was Product Description) --------------------
$ws.Range("A7").Value = "Customer Email"
$ws.Range("B7").Value = "Email:"
$ws.Range("C7").Value = "regex"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""

# --- Row 8: Customer Support (new row) -----------------------------------
$ws.Range("A8").Value = "Customer Support"
$ws.Range("B8").Value = "Phone:"
$ws.Range("C8").Value = "regex"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""

# --- Row 9: Payment Due Date (new row) ------------------------------------
$ws.Range("A9").Value = "Payment Due Date"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "nlp"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "Extract the payment due date or deadline for payment"
